# Adding grid/inline select appearance properties
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "survey" sheet: insert a new "appearance" column (C), fill it in for
#    the birds grid select and the new inline yes_no screen, and add the
#    new begin screen / select_one yes_no (x3) / end screen block in the
#    rows that were already reserved (9-13).
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Insert a blank column before the old "condition" column (C), shifting
# condition/name/label from C/D/E to D/E/F.
$survey.Columns("C:C").Insert()

# New "appearance" header
$survey.Cells.Item(1, 3).Value = "appearance"

# "grid" appearance on the image based "select_one birds" row
$survey.Cells.Item(5, 3).Value = "grid"

# New block: begin screen / select_one yes_no x3 (inline) / end screen
$survey.Cells.Item(9, 2).Value = "begin screen"

$survey.Cells.Item(10, 2).Value = "select_one yes_no"
$survey.Cells.Item(10, 3).Value = "inline"
$survey.Cells.Item(10, 5).Value = "i1"
$survey.Cells.Item(10, 6).Value = "Choose one:"

$survey.Cells.Item(11, 2).Value = "select_one yes_no"
$survey.Cells.Item(11, 3).Value = "inline"
$survey.Cells.Item(11, 5).Value = "i2"
$survey.Cells.Item(11, 6).Value = "Choose one:"

$survey.Cells.Item(12, 2).Value = "select_one yes_no"
$survey.Cells.Item(12, 3).Value = "inline"
$survey.Cells.Item(12, 5).Value = "i3"
$survey.Cells.Item(12, 6).Value = "Choose one:"

$survey.Cells.Item(13, 2).Value = "end screen"

# ---------------------------------------------------------------------
# 2) "choices" sheet: add a "label" column (D) and a new yes_no choice
#    list (yes/Yes, no/No).
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Cells.Item(1, 4).Value = "label"

$choices.Cells.Item(12, 1).Value = "yes_no"
$choices.Cells.Item(12, 2).Value = "yes"
$choices.Cells.Item(12, 4).Value = "Yes"

$choices.Cells.Item(13, 1).Value = "yes_no"
$choices.Cells.Item(13, 2).Value = "no"
$choices.Cells.Item(13, 4).Value = "No"

# ---------------------------------------------------------------------
# 3) "queries" sheet: drop the now-unused param.format/param.q columns,
#    widen the uri column, and switch odk_values to a content:// uri.
# ---------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")

$queries.Range("D1:E1").ClearContents()
$queries.Columns("B:B").ColumnWidth = 70.45

$queries.Cells.Item(5, 2).Value = '"content://com.opendatakit.tables.ContentProvider/database_id/table_id/row_id"'
